$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "24.647.18"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "1.686.38"
$ws.Range("E3").Value = "  -0.97%  "
Set-TextValue $ws.Range("D4") "1.003"
$ws.Range("E4").Value = "  +0.61%  "
Set-TextValue $ws.Range("D5") "315.54"
$ws.Range("E5").Value = "  +0.47%  "
Set-TextValue $ws.Range("D6") "1.003"
$ws.Range("E6").Value = "  +0.55%  "
Set-TextValue $ws.Range("D7") "0.3931"
$ws.Range("E7").Value = "  -1.36%  "
Set-TextValue $ws.Range("D8") "0.4034"
$ws.Range("E8").Value = "  -0.74%  "
Set-TextValue $ws.Range("D9") "1.004"
$ws.Range("E9").Value = "  +0.67%  "
Set-TextValue $ws.Range("D10") "1.481"
Set-TextValue $ws.Range("D11") "52.77"
$ws.Range("E11").Value = "  -0.86%  "
Set-TextValue $ws.Range("D12") "0.08803"
$ws.Range("E12").Value = "  -0.04%  "
Set-TextValue $ws.Range("D13") "7.231"
$ws.Range("E13").Value = "  -1.44%  "
Set-TextValue $ws.Range("D14") "23.36"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("E15").Value = "  +7.07%  "
Set-TextValue $ws.Range("D16") "0.00001307"
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("D17").Value = "1.697.21"
$ws.Range("E17").Value = "  -0.22%  "
Set-TextValue $ws.Range("D18") "99.46"
$ws.Range("E18").Value = "  -1.55%  "
Set-TextValue $ws.Range("D19") "0.07010"
$ws.Range("E19").Value = "  -1.53%  "
Set-TextValue $ws.Range("D20") "19.41"
$ws.Range("E20").Value = "  -0.47%  "
Set-TextValue $ws.Range("D21") "6.976"
$ws.Range("E21").Value = "  +3.27%  "
Set-TextValue $ws.Range("D23") "14.26"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").Value = "24.626.61"
$ws.Range("E24").Value = "  -0.16%  "
Set-TextValue $ws.Range("D25") "3.290"
$ws.Range("E25").Value = "  +9.25%  "
$ws.Range("E26").Value = "  +2.32%  "
Set-TextValue $ws.Range("D27") "22.65"
$ws.Range("E27").Value = "  +0.98%  "
Set-TextValue $ws.Range("D28") "162.27"
$ws.Range("E28").Value = "  +1.91%  "
Set-TextValue $ws.Range("D29") "5.180"
$ws.Range("E29").Value = "  +1.12%  "
Set-TextValue $ws.Range("D30") "135.10"
$ws.Range("E30").Value = "  +1.13%  "
Set-TextValue $ws.Range("D31") "7.544"
$ws.Range("E31").Value = "  +1.04%  "
$ws.Range("D32").Value = "1.879.97"
$ws.Range("E32").Value = "  -0.40%  "
Set-TextValue $ws.Range("D33") "1.056"
$ws.Range("E33").Value = "  -3.22%  "
Set-TextValue $ws.Range("D34") "0.08516"
$ws.Range("E34").Value = "  -1.98%  "
$ws.Range("E35").Value = "  -3.57%  "
$ws.Range("E36").Value = "  -0.01%  "
Set-TextValue $ws.Range("D37") "0.2719"
$ws.Range("E37").Value = "  -0.59%  "
Set-TextValue $ws.Range("D38") "1.875"
$ws.Range("E38").Value = "  -4.19%  "
$ws.Range("E39").Value = "  -3.18%  "
Set-TextValue $ws.Range("D40") "0.09147"
$ws.Range("E40").Value = "  +1.49%  "
Set-TextValue $ws.Range("D41") "0.02702"
$ws.Range("E41").Value = "  -2.58%  "
Set-TextValue $ws.Range("D42") "1.461"
$ws.Range("E42").Value = "  -1.23%  "
Set-TextValue $ws.Range("D43") "0.7576"
$ws.Range("E43").Value = "  -1.12%  "
$ws.Range("E44").Value = "  +2.27%  "
Set-TextValue $ws.Range("D45") "2.583"
$ws.Range("E45").Value = "  +4.99%  "
$ws.Range("E46").Value = "  -1.80%  "
Set-TextValue $ws.Range("D47") "4.218"
$ws.Range("E47").Value = "  +1.26%  "
$ws.Range("E48").Value = "  +0.58%  "
Set-TextValue $ws.Range("D49") "139.43"
Set-TextValue $ws.Range("D50") "1.310"
$ws.Range("E50").Value = "  -0.78%  "
Set-TextValue $ws.Range("D51") "0.07961"
$ws.Range("E51").Value = "  -0.60%  "
